# Final Commit of the Day 2/2
# Updates ObjectID keywords on the TestSteps sheet, updates the SearchText
# value on the TestData sheet, and switches the active sheet/selection.

$wb = $excel.ActiveWorkbook
$wsSteps = $wb.Worksheets.Item("TestSteps")
$wsData  = $wb.Worksheets.Item("TestData")

# --- TestSteps (sheet1): rename a few ObjectID values ---
$wsSteps.Cells.Item(2, 2).Value = "txt_search"
$wsSteps.Cells.Item(3, 2).Value = "txt_search"
$wsSteps.Cells.Item(7, 2).Value = "link_agencyObjectProfile_wait"
$wsSteps.Cells.Item(10, 2).Value = "button_save"
$wsSteps.Cells.Item(11, 2).Value = "button_save"

# --- TestData (sheet2): SearchText value becomes numeric ---
$wsData.Cells.Item(2, 4).Value = 9001

# --- View state: TestData becomes the active/selected tab ---
[void]$wsSteps.Range("B15").Select()
[void]$wsData.Activate()
[void]$wsData.Range("D2").Select()
